$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest daily "cotações" (fund quotes) row - 2025-09-11 (serial 45911) -
# mirroring the layout/format of the existing rows above it.
$ws.Range("A7").Value = 45911
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B7").Value = "21,0689"
$ws.Range("C7").Value = "14,9289"
$ws.Range("D7").Value = "14,8044"
$ws.Range("E7").Value = "14,8044"
